$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.401172666666667
$ws.Range("H2").Value = 10.203518
$ws.Range("I2").Value = 0.5101677883321656
$ws.Range("J2").Value = 0.5101677883321655
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.632955
$ws.Range("N2").Value = 94.898865
$ws.Range("O2").Value = 0.5000340016649593
$ws.Range("P2").Value = 0.5000340016649593
$ws.Range("Q2").Value = 107.5891419118967
$ws.Range("R2").Value = 968.3022772070701
$ws.Range("S2").Value = 0.2551012407202947
$ws.Range("T2").Value = 0.2551012407202946

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.401172666666667
$ws.Range("H3").Value = 10.203518
$ws.Range("I3").Value = 0.5101677883321656
$ws.Range("J3").Value = 0.5101677883321655
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.622273333333332
$ws.Range("N3").Value = 25.86682
$ws.Range("O3").Value = 0.1362955132808722
$ws.Range("P3").Value = 0.1362955132808722
$ws.Range("Q3").Value = 29.32584038586222
$ws.Range("R3").Value = 263.93256347276
$ws.Range("S3").Value = 0.06953358057009988
$ws.Range("T3").Value = 0.06953358057009987

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.401172666666667
$ws.Range("H4").Value = 10.203518
$ws.Range("I4").Value = 0.5101677883321656
$ws.Range("J4").Value = 0.5101677883321655
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.70876033333333
$ws.Range("N4").Value = 56.12628100000001
$ws.Range("O4").Value = 0.2957364019791172
$ws.Range("P4").Value = 0.2957364019791172
$ws.Range("Q4").Value = 63.6317242729509
$ws.Range("R4").Value = 572.6855184565582
$ws.Range("S4").Value = 0.1508751861269985
$ws.Range("T4").Value = 0.1508751861269985

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.401172666666667
$ws.Range("H5").Value = 10.203518
$ws.Range("I5").Value = 0.5101677883321656
$ws.Range("J5").Value = 0.5101677883321655
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.297619333333333
$ws.Range("N5").Value = 12.892858
$ws.Range("O5").Value = 0.06793408307505136
$ws.Range("P5").Value = 0.06793408307505136
$ws.Range("Q5").Value = 14.61694540827156
$ws.Range("R5").Value = 131.552508674444
$ws.Range("S5").Value = 0.03465778091477256
$ws.Range("T5").Value = 0.03465778091477255

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.579868
$ws.Range("H6").Value = 7.739604
$ws.Range("I6").Value = 0.3869740471126509
$ws.Range("J6").Value = 0.3869740471126508
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.632955
$ws.Range("N6").Value = 94.898865
$ws.Range("O6").Value = 0.5000340016649593
$ws.Range("P6").Value = 0.5000340016649593
$ws.Range("Q6").Value = 81.60884834993999
$ws.Range("R6").Value = 734.47963514946
$ws.Range("S6").Value = 0.1935001813182233
$ws.Range("T6").Value = 0.1935001813182233

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.579868
$ws.Range("H7").Value = 7.739604
$ws.Range("I7").Value = 0.3869740471126509
$ws.Range("J7").Value = 0.3869740471126508
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.622273333333332
$ws.Range("N7").Value = 25.86682
$ws.Range("O7").Value = 0.1362955132808722
$ws.Range("P7").Value = 0.1362955132808722
$ws.Range("Q7").Value = 22.24432705992
$ws.Range("R7").Value = 200.19894353928
$ws.Range("S7").Value = 0.05274282637759518
$ws.Range("T7").Value = 0.05274282637759517

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.579868
$ws.Range("H8").Value = 7.739604
$ws.Range("I8").Value = 0.3869740471126509
$ws.Range("J8").Value = 0.3869740471126508
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.70876033333333
$ws.Range("N8").Value = 56.12628100000001
$ws.Range("O8").Value = 0.2957364019791172
$ws.Range("P8").Value = 0.2957364019791172
$ws.Range("Q8").Value = 48.266132103636
$ws.Range("R8").Value = 434.3951889327241
$ws.Range("S8").Value = 0.1144423123523927
$ws.Range("T8").Value = 0.1144423123523927

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.579868
$ws.Range("H9").Value = 7.739604
$ws.Range("I9").Value = 0.3869740471126509
$ws.Range("J9").Value = 0.3869740471126508
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.297619333333333
$ws.Range("N9").Value = 12.892858
$ws.Range("O9").Value = 0.06793408307505136
$ws.Range("P9").Value = 0.06793408307505136
$ws.Range("Q9").Value = 11.087290594248
$ws.Range("R9").Value = 99.785615348232
$ws.Range("S9").Value = 0.02628872706443966
$ws.Range("T9").Value = 0.02628872706443966

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.6857320000000001
$ws.Range("H10").Value = 2.057196
$ws.Range("I10").Value = 0.1028581645551836
$ws.Range("J10").Value = 0.1028581645551836
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.632955
$ws.Range("N10").Value = 94.898865
$ws.Range("O10").Value = 0.5000340016649593
$ws.Range("P10").Value = 0.5000340016649593
$ws.Range("Q10").Value = 21.69172949806
$ws.Range("R10").Value = 195.22556548254
$ws.Range("S10").Value = 0.05143257962644133
$ws.Range("T10").Value = 0.05143257962644131

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.6857320000000001
$ws.Range("H11").Value = 2.057196
$ws.Range("I11").Value = 0.1028581645551836
$ws.Range("J11").Value = 0.1028581645551836
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.622273333333332
$ws.Range("N11").Value = 25.86682
$ws.Range("O11").Value = 0.1362955132808722
$ws.Range("P11").Value = 0.1362955132808722
$ws.Range("Q11").Value = 5.912568737413333
$ws.Range("R11").Value = 53.21311863672
$ws.Range("S11").Value = 0.01401910633317717
$ws.Range("T11").Value = 0.01401910633317716

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.6857320000000001
$ws.Range("H12").Value = 2.057196
$ws.Range("I12").Value = 0.1028581645551836
$ws.Range("J12").Value = 0.1028581645551836
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.70876033333333
$ws.Range("N12").Value = 56.12628100000001
$ws.Range("O12").Value = 0.2957364019791172
$ws.Range("P12").Value = 0.2957364019791172
$ws.Range("Q12").Value = 12.82919564089734
$ws.Range("R12").Value = 115.462760768076
$ws.Range("S12").Value = 0.03041890349972596
$ws.Range("T12").Value = 0.03041890349972595

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.6857320000000001
$ws.Range("H13").Value = 2.057196
$ws.Range("I13").Value = 0.1028581645551836
$ws.Range("J13").Value = 0.1028581645551836
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.297619333333333
$ws.Range("N13").Value = 12.892858
$ws.Range("O13").Value = 0.06793408307505136
$ws.Range("P13").Value = 0.06793408307505136
$ws.Range("Q13").Value = 2.947015100685334
$ws.Range("R13").Value = 26.523135906168
$ws.Range("S13").Value = 0.006987575095839146
$ws.Range("T13").Value = 0.006987575095839144
